$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("完成下周工作计划"): push the due date out ~8 days and mark it in progress
$ws.Range("B2").Value2 = 46052.6666666667
$ws.Range("E2").Value = "进行中"

# Row 3 ("牙医检查"): due date moved a couple of days later
$ws.Range("B3").Value2 = 46052

# Row 4 ("提交健身房会员申请"): due date pushed a month out
$ws.Range("B4").Value2 = 46077

# Move the active cell/selection (cosmetic, matches the saved view state)
$null = $ws.Range("F15").Select()
